# Update symbol list (prices / 1h volume %) for cryptos.xlsx
# Values are forced to text (leading apostrophe) so Excel keeps them as
# literal strings instead of re-interpreting them as numbers/percentages,
# matching the original inline-string cell layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'323.34"
$ws.Range("E2").Value = "'-2.03%"

$ws.Range("D3").Value = "'39.80"
$ws.Range("E3").Value = "'-1.64%"

$ws.Range("D4").Value = "'5.885"
$ws.Range("E4").Value = "'11.54%"

$ws.Range("D5").Value = "'0.08028"
$ws.Range("E5").Value = "'-0.85%"

$ws.Range("D6").Value = "'4.590"
$ws.Range("E6").Value = "'1.48%"

$ws.Range("D7").Value = "'8.667"
$ws.Range("E7").Value = "'0.09%"

$ws.Range("D8").Value = "'1.926"
$ws.Range("E8").Value = "'0.72%"

$ws.Range("D9").Value = "'0.9330"
$ws.Range("E9").Value = "'-0.42%"

$ws.Range("D10").Value = "'0.1272"
$ws.Range("E10").Value = "'-8.25%"

$ws.Range("D11").Value = "'0.1972"
$ws.Range("E11").Value = "'0.39%"

$ws.Range("E12").Value = "'20.02%"

$ws.Range("D13").Value = "'0.09122"
$ws.Range("E13").Value = "'-0.71%"

$ws.Range("D14").Value = "'0.03538"
$ws.Range("E14").Value = "'3.15%"

$ws.Range("D15").Value = "'0.1052"
$ws.Range("E15").Value = "'9.96%"

$ws.Range("D16").Value = "'0.001292"
$ws.Range("E16").Value = "'-7.30%"

$ws.Range("D17").Value = "'0.006119"
$ws.Range("E17").Value = "'3.48%"

$ws.Range("E18").Value = "'-0.33%"

$ws.Range("E19").Value = "'-0.51%"

$ws.Range("D20").Value = "'0.3564"
$ws.Range("E20").Value = "'1.15%"

$ws.Range("D21").Value = "'0.1419"
$ws.Range("E21").Value = "'8.06%"

$ws.Range("D22").Value = "'0.2409"
$ws.Range("E22").Value = "'-6.17%"

$ws.Range("D23").Value = "'0.04411"
$ws.Range("E23").Value = "'-0.76%"

$ws.Range("E24").Value = "'3.16%"

$ws.Range("E25").Value = "'0.38%"

$ws.Range("D26").Value = "'0.0001140"
$ws.Range("E26").Value = "'-11.66%"

$ws.Range("D39").Value = "'0.02439"
$ws.Range("E39").Value = "'-2.84%"

$ws.Range("D40").Value = "'0.05247"
$ws.Range("E40").Value = "'0.45%"

$ws.Range("D41").Value = "'0.007407"
$ws.Range("E41").Value = "'-3.57%"

$ws.Range("D42").Value = "'0.009478"
$ws.Range("E42").Value = "'5.77%"

$ws.Range("E43").Value = "'-1.69%"

$ws.Range("D44").Value = "'0.002120"
$ws.Range("E44").Value = "'-2.34%"

$ws.Range("D45").Value = "'0.009971"
$ws.Range("E45").Value = "'10.87%"

$ws.Range("D46").Value = "'0.00006739"
$ws.Range("E46").Value = "'1.65%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.06%"

$ws.Range("D48").Value = "'0.003000"
$ws.Range("E48").Value = "'-10.19%"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.06%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.06%"
